$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# Update the text in C8: the use case step text changed from the painter
# "asking" to see the services performed to "indicating" that he wants to.
$ws.Range("C8").Value = "Indica que pretende ver Serviços efectuados"

# D11 keeps the same visible text ("Apresenta os Serviços efectuados nessa
# data") - only the shared-string ordering changed in the underlying XML,
# so no visible value change is required here, but set it explicitly to be
# safe/idempotent.
$ws.Range("D11").Value = "Apresenta os Serviços efectuados nessa data"

# Move the active selection from D10 to C9.
$ws.Range("C9").Select()
